$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2300
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 2450
$ws.Range("K31").Value = 6000
$ws.Range("L31").Value = 7350
$ws.Range("M31").Value = -5770
$ws.Range("N31").Value = -7810
$ws.Range("H33").Value = 101.181816
$ws.Range("I33").Value = 98.96666999999999
$ws.Range("K33").Value = 98.96666999999999
$ws.Range("M33").Value = 130.03333
$ws.Range("H51").Value = 2669.4783
$ws.Range("I51").Value = 1944.4445
$ws.Range("J51").Value = 3135.5715
$ws.Range("K51").Value = 1944.4445
$ws.Range("L51").Value = 3135.5715
$ws.Range("M51").Value = -1460.4445
$ws.Range("N51").Value = -4103.5715
$ws.Range("H125").Value = 900
$ws.Range("I125").Value = 716.6667
$ws.Range("K125").Value = 6450.0003
$ws.Range("M125").Value = -3990.0003
$ws.Range("H132").Value = 1521.5862
$ws.Range("I132").Value = 1513.6727
$ws.Range("J132").Value = 1666.6666
$ws.Range("K132").Value = 4541.0181
$ws.Range("L132").Value = 4999.9998
$ws.Range("M132").Value = -2011.0181
$ws.Range("N132").Value = -10059.9998
$ws.Range("H137").Value = 825.325
$ws.Range("I137").Value = 737.4706
$ws.Range("J137").Value = 890.26086
$ws.Range("K137").Value = 2212.4118
$ws.Range("L137").Value = 2670.78258
$ws.Range("M137").Value = 337.5882000000001
$ws.Range("N137").Value = -7770.78258
$ws.Range("H138").Value = 1233.46
$ws.Range("I138").Value = 657.1475
$ws.Range("J138").Value = 2134.8718
$ws.Range("K138").Value = 1971.4425
$ws.Range("L138").Value = 6404.6154
$ws.Range("M138").Value = 3168.5575
$ws.Range("N138").Value = -16684.6154
$ws.Range("H141").Value = 2438.7673
$ws.Range("I141").Value = 831.08826
$ws.Range("K141").Value = 2493.26478
$ws.Range("M141").Value = 2686.73522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6298.17
$ws.Range("I32").Value = 5369.2856
$ws.Range("J32").Value = 11174.8125
$ws.Range("K32").Value = 5369.2856
$ws.Range("L32").Value = 11174.8125
$ws.Range("M32").Value = -5082.2856
$ws.Range("N32").Value = -11748.8125
$ws.Range("H45").Value = 1791.6666
$ws.Range("I45").Value = 1075.7273
$ws.Range("J45").Value = 2579.2
$ws.Range("K45").Value = 1075.7273
$ws.Range("L45").Value = 2579.2
$ws.Range("M45").Value = -698.7273
$ws.Range("N45").Value = -3333.2
$ws.Range("H61").Value = 768.12726
$ws.Range("I61").Value = 674.1957
$ws.Range("J61").Value = 1248.2222
$ws.Range("K61").Value = 674.1957
$ws.Range("L61").Value = 1248.2222
$ws.Range("M61").Value = -462.1957
$ws.Range("N61").Value = -1672.2222
$ws.Range("H132").Value = 963.78845
$ws.Range("I132").Value = 840
$ws.Range("J132").Value = 1425.1818
$ws.Range("K132").Value = 2520
$ws.Range("L132").Value = 4275.5454
$ws.Range("M132").Value = 10
$ws.Range("N132").Value = -9335.545399999999
$ws.Range("H136").Value = 768.12726
$ws.Range("I136").Value = 674.1957
$ws.Range("J136").Value = 1248.2222
$ws.Range("K136").Value = 2022.5871
$ws.Range("L136").Value = 3744.6666
$ws.Range("M136").Value = 527.4129
$ws.Range("N136").Value = -8844.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 26125
$ws.Range("J100").Value = 26125
$ws.Range("L100").Value = 26125
$ws.Range("N100").Value = -28289
$ws.Range("H134").Value = 13476
$ws.Range("I134").Value = 1253.0845
$ws.Range("J134").Value = 85794.914
$ws.Range("K134").Value = 3759.2535
$ws.Range("L134").Value = 257384.742
$ws.Range("M134").Value = -1224.2535
$ws.Range("N134").Value = -262454.742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3380.276
$ws.Range("I31").Value = 3275.1304
$ws.Range("J31").Value = 3783.3333
$ws.Range("K31").Value = 3275.1304
$ws.Range("L31").Value = 3783.3333
$ws.Range("M31").Value = -2980.1304
$ws.Range("N31").Value = -4373.3333
$ws.Range("H34").Value = 3380.276
$ws.Range("I34").Value = 3275.1304
$ws.Range("J34").Value = 3783.3333
$ws.Range("K34").Value = 3275.1304
$ws.Range("L34").Value = 3783.3333
$ws.Range("M34").Value = -3073.1304
$ws.Range("N34").Value = -4187.3333
$ws.Range("H58").Value = 2991.689
$ws.Range("I58").Value = 735.4400000000001
$ws.Range("J58").Value = 5812
$ws.Range("K58").Value = 735.4400000000001
$ws.Range("L58").Value = 5812
$ws.Range("M58").Value = -532.4400000000001
$ws.Range("N58").Value = -6218
$ws.Range("H97").Value = 22000
$ws.Range("J97").Value = 22000
$ws.Range("L97").Value = 22000
$ws.Range("N97").Value = -23982
$ws.Range("H132").Value = 1589.9868
$ws.Range("I132").Value = 1006.5227
$ws.Range("J132").Value = 2392.25
$ws.Range("K132").Value = 3019.5681
$ws.Range("L132").Value = 7176.75
$ws.Range("M132").Value = -489.5681
$ws.Range("N132").Value = -12236.75
$ws.Range("H134").Value = 1040.3066
$ws.Range("I134").Value = 1020.88464
$ws.Range("J134").Value = 1084.2174
$ws.Range("K134").Value = 3062.65392
$ws.Range("L134").Value = 3252.6522
$ws.Range("M134").Value = -527.6539199999997
$ws.Range("N134").Value = -8322.6522
$ws.Range("H136").Value = 2991.689
$ws.Range("I136").Value = 735.4400000000001
$ws.Range("J136").Value = 5812
$ws.Range("K136").Value = 2206.32
$ws.Range("L136").Value = 17436
$ws.Range("M136").Value = 343.6799999999998
$ws.Range("N136").Value = -22536

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 592.0714
$ws.Range("I5").Value = 558.7027
$ws.Range("J5").Value = 839
$ws.Range("K5").Value = 1676.1081
$ws.Range("L5").Value = 2517
$ws.Range("M5").Value = -1564.1081
$ws.Range("N5").Value = -2741
$ws.Range("H45").Value = 1081.6
$ws.Range("J45").Value = 1238
$ws.Range("L45").Value = 3714
$ws.Range("N45").Value = -4778
$ws.Range("H113").Value = 622.8570999999999
$ws.Range("I113").Value = 707.2727
$ws.Range("J113").Value = 568.2353000000001
$ws.Range("K113").Value = 2121.8181
$ws.Range("L113").Value = 1704.7059
$ws.Range("M113").Value = 48.18190000000004
$ws.Range("N113").Value = -6044.7059
$ws.Range("H135").Value = 592.0714
$ws.Range("I135").Value = 558.7027
$ws.Range("J135").Value = 839
$ws.Range("K135").Value = 5028.3243
$ws.Range("L135").Value = 7551
$ws.Range("M135").Value = -2493.3243
$ws.Range("N135").Value = -12621

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4133.3335
$ws.Range("I70").Value = 4003.5715
$ws.Range("K70").Value = 4003.5715
$ws.Range("M70").Value = -3733.5715
$ws.Range("H73").Value = 4133.3335
$ws.Range("I73").Value = 4003.5715
$ws.Range("K73").Value = 4003.5715
$ws.Range("M73").Value = -3067.5715
$ws.Range("H82").Value = 31833
$ws.Range("J82").Value = 31833
$ws.Range("L82").Value = 31833
$ws.Range("N82").Value = -32599
$ws.Range("H85").Value = 31833
$ws.Range("J85").Value = 31833
$ws.Range("L85").Value = 31833
$ws.Range("N85").Value = -34485
$ws.Range("H132").Value = 1783.8254
$ws.Range("I132").Value = 1766.7142
$ws.Range("J132").Value = 1805.2142
$ws.Range("K132").Value = 5300.142599999999
$ws.Range("L132").Value = 5415.642599999999
$ws.Range("M132").Value = -2770.142599999999
$ws.Range("N132").Value = -10475.6426
$ws.Range("H140").Value = 70960
$ws.Range("J140").Value = 70960
$ws.Range("L140").Value = 70960
$ws.Range("N140").Value = -81320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5425.6924
$ws.Range("I61").Value = 6503.4
$ws.Range("K61").Value = 6503.4
$ws.Range("M61").Value = -6301.4
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 5425.6924
$ws.Range("I113").Value = 6503.4
$ws.Range("K113").Value = 6503.4
$ws.Range("M113").Value = -4333.4
$ws.Range("H132").Value = 1652.8088
$ws.Range("I132").Value = 1451.6111
$ws.Range("J132").Value = 2428.8572
$ws.Range("K132").Value = 4354.8333
$ws.Range("L132").Value = 7286.571599999999
$ws.Range("M132").Value = -1824.8333
$ws.Range("N132").Value = -12346.5716
$ws.Range("H136").Value = 1712.3103
$ws.Range("I136").Value = 959.2558
$ws.Range("J136").Value = 3871.0667
$ws.Range("K136").Value = 2877.7674
$ws.Range("L136").Value = 11613.2001
$ws.Range("M136").Value = -327.7674000000002
$ws.Range("N136").Value = -16713.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1272.6364
$ws.Range("I122").Value = 1149.25
$ws.Range("J122").Value = 1601.6666
$ws.Range("K122").Value = 3447.75
$ws.Range("L122").Value = 4804.9998
$ws.Range("M122").Value = -997.75
$ws.Range("N122").Value = -9704.9998
$ws.Range("H126").Value = 617.875
$ws.Range("I126").Value = 617.875
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1853.625
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 616.375
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 718.92065
$ws.Range("I132").Value = 630.0682
$ws.Range("J132").Value = 924.6842
$ws.Range("K132").Value = 1890.2046
$ws.Range("L132").Value = 2774.0526
$ws.Range("M132").Value = 639.7954
$ws.Range("N132").Value = -7834.0526
$ws.Range("H136").Value = 890.1951
$ws.Range("I136").Value = 963.7931
$ws.Range("J136").Value = 712.3333
$ws.Range("K136").Value = 2891.3793
$ws.Range("L136").Value = 2136.9999
$ws.Range("M136").Value = -341.3793000000001
$ws.Range("N136").Value = -7236.9999
